$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.905.75"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "1.809.36"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.47"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4615"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07376"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8735"
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "1.807.03"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.357"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.530"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.82"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07049"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008689"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "26.904.25"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.323"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").Value = "2.023.97"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.20"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -6.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.299"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.89"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08901"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7529"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.155"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.913"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.444"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.101"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01970"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05232"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.423"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5289"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.172"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1661"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4968"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.33"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.94"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.670"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("E51").Value = "  -1.51%  "
